$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:D1) stays the same: Brand, Model, images_name, Image_Type_Nmae

# Update row 2 -> Neptune / PS-50 Power product
$ws.Range("A2").Value = "Neptune"
$ws.Range("B2").Value = "PS-50 Power"
$ws.Range("C2").Value = "['PS-50Powerimg0-ps-50-power-56-1650083453.png', 'PS-50Powerimg1-default-image.png', 'PS-50Powerimg2-ps-50-power-56-1650083453.png']"
$ws.Range("D2").Value = "product"

# Update row 3 -> Neptune / NF-10B Manual product
$ws.Range("A3").Value = "Neptune"
$ws.Range("B3").Value = "NF-10B Manual"
$ws.Range("C3").Value = "['NF-10BManualimg0-nf-10b-manual-56-1650083150.png', 'NF-10BManualimg1-default-image.png', 'NF-10BManualimg2-nf-10b-manual-56-1650083150.png']"
$ws.Range("D3").Value = "product"

# Remove the old fourth row entirely (table shrinks from 3 products to 2)
$ws.Rows.Item(4).Delete()
